$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Full1")

# Update the "SI"/"NO" completion markers in rows 10 and 11
$ws.Range("D10").Value = "SI"
$ws.Range("F10").Value = "SI"
$ws.Range("I10").Value = "SI"
$ws.Range("F11").Value = "SI"
$ws.Range("I11").Value = "SI"

# Reflect the active selection recorded in the saved workbook
$ws.Range("H10").Select()

$wb.Save()
